# Update "想去人数" (F) and occasionally "最低票价" (G) figures across the
# three sheets that share the same conference data (展览, 演出, 全部类型).
# Values below were taken from the target commit's XML diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 34
$ws1.Range("F4").Value  = 409
$ws1.Range("F5").Value  = 5378
$ws1.Range("F6").Value  = 5378
$ws1.Range("F7").Value  = 183
$ws1.Range("F9").Value  = 537
$ws1.Range("F11").Value = 1205
$ws1.Range("F12").Value = 6333
$ws1.Range("F13").Value = 35
$ws1.Range("F15").Value = 110
$ws1.Range("F16").Value = 3242
$ws1.Range("F17").Value = 267
$ws1.Range("F19").Value = 257
$ws1.Range("F20").Value = 4044
$ws1.Range("F24").Value = 3961
$ws1.Range("F25").Value = 192
$ws1.Range("F26").Value = 191
$ws1.Range("F28").Value = 249
$ws1.Range("F32").Value = 129
$ws1.Range("F36").Value = 30
$ws1.Range("F37").Value = 7069
$ws1.Range("F39").Value = 1161
$ws1.Range("F40").Value = 564
$ws1.Range("F43").Value = 1427
$ws1.Range("F45").Value = 757
$ws1.Range("F46").Value = 3018
$ws1.Range("F47").Value = 324
$ws1.Range("F49").Value = 794
$ws1.Range("F50").Value = 984

# --- Sheet: 演出 -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 93
$ws2.Range("G16").Value = 128
$ws2.Range("F22").Value = 57
$ws2.Range("F25").Value = 835

# --- Sheet: 全部类型 ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 34
$ws4.Range("F7").Value  = 409
$ws4.Range("F8").Value  = 5378
$ws4.Range("F9").Value  = 5378
$ws4.Range("F10").Value = 183
$ws4.Range("F12").Value = 93
$ws4.Range("F13").Value = 537
$ws4.Range("F15").Value = 1205
$ws4.Range("F16").Value = 6334
$ws4.Range("F17").Value = 35
$ws4.Range("F19").Value = 110
$ws4.Range("F20").Value = 3242
$ws4.Range("F21").Value = 267
$ws4.Range("F23").Value = 257
$ws4.Range("F24").Value = 4044
$ws4.Range("F25").Value = 3961
$ws4.Range("F26").Value = 192
$ws4.Range("F27").Value = 191
$ws4.Range("F28").Value = 249
$ws4.Range("F32").Value = 129
$ws4.Range("F36").Value = 7069
$ws4.Range("F38").Value = 1161
$ws4.Range("F39").Value = 564
$ws4.Range("F43").Value = 1427
$ws4.Range("F45").Value = 757
$ws4.Range("F46").Value = 3019
$ws4.Range("F47").Value = 324
$ws4.Range("F48").Value = 794
$ws4.Range("F49").Value = 984
